# Weekly price-list update: a new "Haba" price record (week of 2023-11-16,
# serial date 45246) is inserted as the new row 26 of the data table,
# pushing every following record down by one row (old row 26 becomes row
# 27, ..., old row 98 becomes row 99). The sheet's used range grows from
# A1:R98 to A1:R99.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 26..98 down by inserting a brand-new blank row at position 26.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new record's data.
$ws.Cells.Item(26, 1).Value  = 7
$ws.Cells.Item(26, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(26, 3).Value  = "Ñuble"
$ws.Cells.Item(26, 4).Value  = 45246
$ws.Cells.Item(26, 5).Value  = 16
$ws.Cells.Item(26, 6).Value  = 100112026
$ws.Cells.Item(26, 7).Value  = "Haba"
$ws.Cells.Item(26, 8).Value  = "Sin especificar"
$ws.Cells.Item(26, 9).Value  = "Primera"
$ws.Cells.Item(26, 10).Value = 100
$ws.Cells.Item(26, 11).Value = 10000
$ws.Cells.Item(26, 12).Value = 10000
$ws.Cells.Item(26, 13).Value = 10000
$ws.Cells.Item(26, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(26, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(26, 16).Value = 400
$ws.Cells.Item(26, 17).Value = 25
$ws.Cells.Item(26, 18).Value = "Hortaliza"
